$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (attendance) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 332
$ws1.Range("F4").Value = 1367
$ws1.Range("F5").Value = 656

# Sheet "全部类型" (All Types) - update "想去人数" (attendance) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 332
$ws4.Range("F4").Value = 1367
$ws4.Range("F6").Value = 656
